$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A9 with the new feedback text (within the shared string at index 7)
$ws.Range("A9").Value = "b'Moro!' @ 10.5.2017, 20:54"

# Delete rows 10 to 16 (old rows beyond the max length check)
$ws.Range("A10:A16").EntireRow.Delete()

# Set column A width (closest achievable value to the target stored width of 98.6640625)
$ws.Columns.Item(1).ColumnWidth = 97.83

# Update selection
$ws.Range("C12").Select()
